# Apply the "PO Forecast" update to the workbook.
$wb = $excel.ActiveWorkbook

# --- 1. Rename header labels on the existing sheets -----------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after "Monthly Trend" -------------
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Match the header / date-column formatting used on the other sheets by
# copying it over before the values are written (PasteSpecial only touches
# formatting here, but it also clears the destination's contents, so this
# must happen before the Value assignments below).
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = @(
    @(45207.99999999999, 7,  -15.06723053151846, 32.83305662570882),
    @(45228.99999999999, 8,  -14.85501266316585, 32.39803002944067),
    @(45298.99999999999, 10, -15.27434866013709, 34.5548044033676),
    @(45305.99999999999, 10, -12.12371817502301, 34.63586376288541),
    @(45319.99999999999, 11, -13.41606269662277, 32.44477269692267),
    @(45333.99999999999, 11, -12.96895860367505, 34.1597702491769),
    @(45424.99999999999, 14, -7.82719259038907,  38.56540241108802),
    @(45431.99999999999, 14, -8.989458802677321, 37.01802385457864),
    @(45438.99999999999, 15, -8.704476636794602, 37.41601127941094),
    @(45445.99999999999, 15, -9.002445849181031, 39.71469420578054),
    @(45452.99999999999, 15, -8.21071275241006,  38.98609319410157),
    @(45459.99999999999, 15, -10.29489303369836, 38.43128734481732),
    @(45466.99999999999, 16, -8.474917430972821, 38.31050932684986),
    @(45473.99999999999, 16, -8.805150309716712, 39.72167873009519),
    @(45480.99999999999, 16, -6.852819194776938, 39.55220261023613),
    @(45487.99999999999, 16, -5.80336712251033,  40.36102106483953)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}
